# Update column C ("Förändrad") date value from 2024-05-19 (45431)
# to 2024-05-20 (45432) for rows 2 through 28.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45431) {
        $cell.Value2 = 45432
    }
}
